# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of cell address -> new value (same updates apply to both sheets).
$updates = @{
    "F2"  = 1881
    "F4"  = 111
    "F6"  = 14
    "F7"  = 1557
    "F9"  = 607
    "F10" = 364
    "F12" = 19
    "F17" = 104
    "F18" = 121
    "F19" = 3634
    "F20" = 2
    "F21" = 3
    "F22" = 427
    "F23" = 326
    "F24" = 487
    "F25" = 251
    "F26" = 341
    "F28" = 1412
    "F29" = 139
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
